# Apply the cryptos list refresh described in the commit
# ("Updated cryptos list on Fri Jun 23 17:25:59 UTC 2023 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.075.01"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "1.916.33"
$ws.Range("E3").Value = "  +1.95%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").Value = "'245.63"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").Value = "'0.4981"
$ws.Range("E7").Value = "  +0.99%  "

$ws.Range("E8").Value = "  +3.28%  "

$ws.Range("D9").Value = "'0.06850"
$ws.Range("E9").Value = "  +3.63%  "

$ws.Range("D10").Value = "1.911.62"
$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("D11").Value = "'17.01"
$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").Value = "'0.07317"
$ws.Range("E12").Value = "  +1.74%  "

$ws.Range("E13").Value = "  +5.21%  "

$ws.Range("D14").Value = "'0.6835"
$ws.Range("E14").Value = "  +2.96%  "

$ws.Range("D15").Value = "'5.083"
$ws.Range("E15").Value = "  +5.21%  "

$ws.Range("D16").Value = "31.055.98"
$ws.Range("E16").Value = "  +3.71%  "

$ws.Range("D17").Value = "'0.000008030"
$ws.Range("E17").Value = "  +1.91%  "

$ws.Range("D18").Value = "'13.40"
$ws.Range("E18").Value = "  +5.03%  "

$ws.Range("D19").Value = "'0.9991"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "2.160.84"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").Value = "'0.9954"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "'4.881"
$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23").Value = "'174.56"
$ws.Range("E23").Value = "  +28.00%  "

$ws.Range("D24").Value = "'6.072"
$ws.Range("E24").Value = "  +9.23%  "

$ws.Range("D25").Value = "'9.332"
$ws.Range("E25").Value = "  +2.11%  "

$ws.Range("D26").Value = "'151.68"
$ws.Range("E26").Value = "  +2.74%  "

$ws.Range("D27").Value = "'18.17"
$ws.Range("E27").Value = "  +8.31%  "

$ws.Range("D28").Value = "'1.952"
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("D29").Value = "'1.414"
$ws.Range("E29").Value = "  +2.33%  "

$ws.Range("D30").Value = "'4.376"
$ws.Range("E30").Value = "  +4.77%  "

$ws.Range("D31").Value = "'0.08946"
$ws.Range("E31").Value = "  +3.67%  "

$ws.Range("D32").Value = "'4.068"
$ws.Range("E32").Value = "  +3.26%  "

$ws.Range("D33").Value = "'0.05262"
$ws.Range("E33").Value = "  +6.01%  "

$ws.Range("D34").Value = "'0.7493"
$ws.Range("E34").Value = "  +6.78%  "

$ws.Range("D35").Value = "'1.144"
$ws.Range("E35").Value = "  +3.11%  "

$ws.Range("D36").Value = "'2.662"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").Value = "'0.01917"
$ws.Range("E37").Value = "  +17.05%  "

$ws.Range("D38").Value = "'2.749"
$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("D39").Value = "'2.206"
$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("D40").Value = "'0.9441"
$ws.Range("E40").Value = "  +1.55%  "

$ws.Range("D41").Value = "'5.967"
$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("D42").Value = "'0.4365"
$ws.Range("E42").Value = "  +4.45%  "

$ws.Range("D43").Value = "'105.13"
$ws.Range("E43").Value = "  +3.22%  "

$ws.Range("D44").Value = "'7.857"
$ws.Range("E44").Value = "  +3.71%  "

$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "'0.1329"
$ws.Range("E46").Value = "  +5.83%  "

$ws.Range("D47").Value = "'0.05863"
$ws.Range("E47").Value = "  +2.74%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.608"
$ws.Range("E48").Value = "  +5.26%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.3905"
$ws.Range("E49").Value = "  +5.60%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'33.35"
$ws.Range("E50").Value = "  +2.55%  "

$ws.Range("D51").Value = "'1.391"
$ws.Range("E51").Value = "  +4.46%  "
